# Auto-generated script to apply scheduled-runner price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 390
$ws.Range("I12").Value = 281
$ws.Range("J12").Value = 499
$ws.Range("K12").Value = 281
$ws.Range("L12").Value = 499
$ws.Range("M12").Value = -111
$ws.Range("N12").Value = -839
$ws.Range("H19").Value = 2296
$ws.Range("I19").Value = 2289
$ws.Range("J19").Value = 2299.5
$ws.Range("K19").Value = 2289
$ws.Range("L19").Value = 2299.5
$ws.Range("M19").Value = -2114
$ws.Range("N19").Value = -2649.5
$ws.Range("H28").Value = 2927.5
$ws.Range("J28").Value = 3282.1667
$ws.Range("L28").Value = 3282.1667
$ws.Range("N28").Value = -4252.1667
$ws.Range("H33").Value = 405.9091
$ws.Range("I33").Value = 412.8
$ws.Range("K33").Value = 412.8
$ws.Range("M33").Value = -183.8
$ws.Range("H68").Value = 50000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 50000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 50000
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -51498
$ws.Range("H71").Value = 50000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 50000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 150000
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -157488
$ws.Range("H97").Value = 3903.3333
$ws.Range("J97").Value = 3903.3333
$ws.Range("L97").Value = 11709.9999
$ws.Range("N97").Value = -12701.9999
$ws.Range("H107").Value = 1595.8889
$ws.Range("I107").Value = 1285.6666
$ws.Range("J107").Value = 2216.3333
$ws.Range("K107").Value = 1285.6666
$ws.Range("L107").Value = 2216.3333
$ws.Range("M107").Value = 634.3334
$ws.Range("N107").Value = -6056.3333
$ws.Range("H141").Value = 13865.667
$ws.Range("I141").Value = 13865.667
$ws.Range("K141").Value = 41597.001
$ws.Range("M141").Value = -36417.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4217.6665
$ws.Range("I32").Value = 4217.6665
$ws.Range("K32").Value = 4217.6665
$ws.Range("M32").Value = -3930.6665
$ws.Range("H45").Value = 4346
$ws.Range("I45").Value = 3512
$ws.Range("K45").Value = 3512
$ws.Range("M45").Value = -3135
$ws.Range("H61").Value = 3104.0454
$ws.Range("I61").Value = 2686.4
$ws.Range("K61").Value = 2686.4
$ws.Range("M61").Value = -2474.4
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
$ws.Range("H104").Value = 33333.332
$ws.Range("I104").Value = 40000
$ws.Range("J104").Value = 30000
$ws.Range("K104").Value = 40000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988
$ws.Range("M104").Value = -36506
$ws.Range("H114").Value = 29000
$ws.Range("J114").Value = 29000
$ws.Range("L114").Value = 29000
$ws.Range("N114").Value = -37678
$ws.Range("H132").Value = 3178.7693
$ws.Range("I132").Value = 3376.6667
$ws.Range("K132").Value = 10130.0001
$ws.Range("M132").Value = -7600.000100000001
$ws.Range("H136").Value = 3104.0454
$ws.Range("I136").Value = 2686.4
$ws.Range("K136").Value = 8059.200000000001
$ws.Range("M136").Value = -5509.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1257.7
$ws.Range("I80").Value = 436.75
$ws.Range("J80").Value = 4541.5
$ws.Range("K80").Value = 436.75
$ws.Range("L80").Value = 4541.5
$ws.Range("M80").Value = 561.25
$ws.Range("N80").Value = -6537.5
$ws.Range("H83").Value = 1257.7
$ws.Range("I83").Value = 436.75
$ws.Range("J83").Value = 4541.5
$ws.Range("K83").Value = 2183.75
$ws.Range("L83").Value = 22707.5
$ws.Range("M83").Value = 2808.25
$ws.Range("N83").Value = -32691.5
$ws.Range("H107").Value = 2480
$ws.Range("I107").Value = 2480
$ws.Range("K107").Value = 2480
$ws.Range("M107").Value = -560

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6667898
$ws.Range("J22").Value = 20000696
$ws.Range("L22").Value = 20000696
$ws.Range("N22").Value = -20001396
$ws.Range("H96").Value = 17257.666
$ws.Range("J96").Value = 17257.666
$ws.Range("L96").Value = 17257.666
$ws.Range("N96").Value = -22749.666
$ws.Range("H141").Value = 531250
$ws.Range("J141").Value = 1000000
$ws.Range("L141").Value = 1000000
$ws.Range("N141").Value = -1010360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1981
$ws.Range("I113").Value = 988.5
$ws.Range("K113").Value = 2965.5
$ws.Range("M113").Value = -795.5
$ws.Range("H131").Value = 901
$ws.Range("J131").Value = 1599.5
$ws.Range("L131").Value = 4798.5
$ws.Range("N131").Value = -14878.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 22010
$ws.Range("J27").Value = 22010
$ws.Range("L27").Value = 22010
$ws.Range("N27").Value = -22342

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2199
$ws.Range("I7").Value = 2199
$ws.Range("K7").Value = 2199
$ws.Range("M7").Value = -2087
$ws.Range("H99").Value = 29999
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H101").Value = 17000
$ws.Range("J101").Value = 17000
$ws.Range("L101").Value = 17000
$ws.Range("N101").Value = -23490
$ws.Range("H126").Value = 2199
$ws.Range("I126").Value = 2199
$ws.Range("K126").Value = 6597
$ws.Range("M126").Value = -4127
$ws.Range("H136").Value = 12861.883
$ws.Range("I136").Value = 13359.538
$ws.Range("J136").Value = 11244.5
$ws.Range("K136").Value = 40078.614
$ws.Range("L136").Value = 33733.5
$ws.Range("M136").Value = -37528.614
$ws.Range("N136").Value = -38833.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 52944.25
$ws.Range("I45").Value = 43331
$ws.Range("J45").Value = 58712.2
$ws.Range("K45").Value = 43331
$ws.Range("L45").Value = 58712.2
$ws.Range("M45").Value = -42840
$ws.Range("N45").Value = -59694.2
$ws.Range("H124").Value = 23960.5
$ws.Range("J124").Value = 23960.5
$ws.Range("L124").Value = 23960.5
$ws.Range("N124").Value = -33780.5
$ws.Range("H132").Value = 3436.0715
$ws.Range("I132").Value = 3141.6667
$ws.Range("K132").Value = 9425.000100000001
$ws.Range("M132").Value = -6895.000100000001
